$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Prestacion" column first (keeps existing IAC.* labels, adds the
# new COV.* ones, including a brand-new row for the extra prestacion).
$ws.Range("A2").Value = "IAC.01"
$ws.Range("A3").Value = "IAC.02"
$ws.Range("A4").Value = "IAC.03"
$ws.Range("A5").Value = "COV.16"
$ws.Range("A6").Value = "COV.17"

# Headers.
$ws.Range("A1").Value = "Prestacion"
$ws.Range("B1").Value = "Grupo"

# Every prestacion now rolls up under the "TEST" group (replaces "Nancy"/"Pablo").
$ws.Range("B2").Value = "TEST"
$ws.Range("B3").Value = "TEST"
$ws.Range("B4").Value = "TEST"
$ws.Range("B5").Value = "TEST"
$ws.Range("B6").Value = "TEST"

# Grow the table to cover the new row (and then some, matching the author's resize).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:B15"))

# Match the author's final cell selection.
$ws.Range("E11").Select()
